$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.3628830909729
$ws.Range("B1").Value = 2.935124635696411
$ws.Range("C1").Value = 4.00922155380249
$ws.Range("D1").Value = 3.173510074615479
$ws.Range("E1").Value = 0.7575874328613281
